$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# spring 24 week 11 inputs
$ws.Range("E3").Value = 1.32
$ws.Range("C4").Value = 1.43
$ws.Range("E6").Value = 1.31
$ws.Range("G6").Value = 1.01
$ws.Range("F7").Value = 1.47
